# Combine station 5 and 5b (see commit message) across all four sheets.
# For each sheet, the unique (non-overlapping) values from the "5b" row are
# folded into the "5" row, any overlapping numeric columns are updated to the
# new combined value, and the now-redundant "5b" row is deleted so every row
# below it shifts up.

$wb = $excel.ActiveWorkbook

# --- Sheet 1: "Littorina TBT" ---
# Row 4 = Station "5", Row 5 = Station "5b"
$ws1 = $wb.Worksheets.Item("Littorina TBT")
$ws1.Cells.Item(4, 5).Value = 35.7     # E4
$ws1.Cells.Item(4, 7).Value = 243      # G4
$ws1.Cells.Item(4, 8).Value = 120.35   # H4 (average of 96.7 and 144)
$ws1.Rows.Item(5).Delete()

# --- Sheet 2: "Littorina ISI" ---
# Row 4 = Station "5", Row 5 = Station "5b"
$ws2 = $wb.Worksheets.Item("Littorina ISI")
$ws2.Cells.Item(4, 3).Value = 0.68   # C4
$ws2.Cells.Item(4, 4).Value = 0.4    # D4
$ws2.Cells.Item(4, 5).Value = 0.52   # E4
$ws2.Rows.Item(5).Delete()

# --- Sheet 3: "Littorina sterile" ---
# Row 4 = Station "5", Row 5 = Station "5b"
$ws3 = $wb.Worksheets.Item("Littorina sterile")
$ws3.Cells.Item(4, 3).Value = 0.68   # C4
$ws3.Cells.Item(4, 4).Value = 0.4    # D4
$ws3.Cells.Item(4, 5).Value = 0.52   # E4
$ws3.Rows.Item(5).Delete()

# --- Sheet 4: "Littorina PRL" ---
# Row 3 = Station "5", Row 4 = Station "5b"
$ws4 = $wb.Worksheets.Item("Littorina PRL")
$ws4.Cells.Item(3, 3).Value = 3.11   # C3
$ws4.Cells.Item(3, 4).Value = 6      # D3
$ws4.Cells.Item(3, 5).Value = 5.22   # E3
$ws4.Rows.Item(4).Delete()
